# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that "Administrator" is listed first, followed by the remaining names in
# their original order, wherever the cell currently contains "Administrator".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*Administrator*") {
        $parts = $val -split ",\s*"
        $rest = $parts | Where-Object { $_ -ne "Administrator" }
        $newParts = @("Administrator") + $rest
        $cell.Value2 = ($newParts -join ", ")
    }
}
